$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.873.47'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '1.659.15'
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9994'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3642'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.59%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.58'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3264'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.135'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07076'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9995'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.071'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.57'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.80%  '
$ws.Range("D15").Value = '1.656.82'
$ws.Range("E15").Value = '  -0.79%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.598'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001047'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06589'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '79.07'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.915'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.77'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.64'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.06%  '
$ws.Range("D24").Value = '24.867.67'
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.444'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.455'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '147.83'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.66'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.44%  '
$ws.Range("D29").Value = '1.838.53'
$ws.Range("E29").Value = '  -1.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.203'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '125.31'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.087'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.771'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -10.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08441'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.649'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.27'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.86%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.282'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.71%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.173'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02263'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06068'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.427'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2069'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9992'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5937'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.882'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.99'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5627'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.17'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.949'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06992'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.190'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.60%  '
